$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 368.5
$ws.Range("J28").Value = 696.6667
$ws.Range("L28").Value = 696.6667
$ws.Range("N28").Value = -1666.6667
$ws.Range("H32").Value = 8755
$ws.Range("I32").Value = 6500
$ws.Range("J32").Value = 9506.666999999999
$ws.Range("K32").Value = 6500
$ws.Range("L32").Value = 9506.666999999999
$ws.Range("M32").Value = -6174
$ws.Range("N32").Value = -10158.667
$ws.Range("H86").Value = 8200
$ws.Range("J86").Value = 8000
$ws.Range("L86").Value = 8000
$ws.Range("N86").Value = -10246
$ws.Range("H89").Value = 8200
$ws.Range("J89").Value = 8000
$ws.Range("L89").Value = 40000
$ws.Range("N89").Value = -51232
$ws.Range("H98").Value = 2200
$ws.Range("I98").Value = 2200
$ws.Range("K98").Value = 2200
$ws.Range("M98").Value = -702
$ws.Range("H116").Value = 5358.6
$ws.Range("I116").Value = 3996.5
$ws.Range("K116").Value = 3996.5
$ws.Range("M116").Value = -554.5
$ws.Range("H122").Value = 2200
$ws.Range("I122").Value = 2200
$ws.Range("K122").Value = 6600
$ws.Range("M122").Value = -4150
$ws.Range("H138").Value = 4773.4023
$ws.Range("I138").Value = 2821.5925
$ws.Range("J138").Value = 5651.717
$ws.Range("K138").Value = 8464.7775
$ws.Range("L138").Value = 16955.151
$ws.Range("M138").Value = -3324.7775
$ws.Range("N138").Value = -27235.151

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2001668
$ws.Range("I2").Value = 2192153.2
$ws.Range("K2").Value = 2192153.2
$ws.Range("M2").Value = -2192040.2
$ws.Range("H32").Value = 23486.64
$ws.Range("I32").Value = 23658.13
$ws.Range("K32").Value = 23658.13
$ws.Range("M32").Value = -23371.13
$ws.Range("H116").Value = 2001668
$ws.Range("I116").Value = 2192153.2
$ws.Range("K116").Value = 2192153.2
$ws.Range("M116").Value = -2189859.2
$ws.Range("H132").Value = 5399.968
$ws.Range("I132").Value = 3288.423
$ws.Range("J132").Value = 16380
$ws.Range("K132").Value = 9865.269
$ws.Range("L132").Value = 49140
$ws.Range("M132").Value = -7335.269
$ws.Range("N132").Value = -54200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2001668
$ws.Range("I3").Value = 2192153.2
$ws.Range("K3").Value = 2192153.2
$ws.Range("M3").Value = -2192039.2
$ws.Range("H107").Value = 2266
$ws.Range("I107").Value = 1899
$ws.Range("K107").Value = 1899
$ws.Range("M107").Value = 21

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1681.7693
$ws.Range("J16").Value = 2248
$ws.Range("L16").Value = 2248
$ws.Range("N16").Value = -2822
$ws.Range("H62").Value = 24169.334
$ws.Range("I62").Value = 20005
$ws.Range("J62").Value = 25002.2
$ws.Range("K62").Value = 20005
$ws.Range("L62").Value = 25002.2
$ws.Range("M62").Value = -19381
$ws.Range("N62").Value = -26250.2
$ws.Range("H65").Value = 24169.334
$ws.Range("I65").Value = 20005
$ws.Range("J65").Value = 25002.2
$ws.Range("K65").Value = 100025
$ws.Range("L65").Value = 125011
$ws.Range("M65").Value = -96905
$ws.Range("N65").Value = -131251
$ws.Range("H107").Value = 584.4286
$ws.Range("I107").Value = 584.4286
$ws.Range("K107").Value = 584.4286
$ws.Range("M107").Value = 1335.5714
$ws.Range("H113").Value = 1681.7693
$ws.Range("J113").Value = 2248
$ws.Range("L113").Value = 2248
$ws.Range("N113").Value = -6588
$ws.Range("H141").Value = 192753.38
$ws.Range("J141").Value = 214930.14
$ws.Range("L141").Value = 214930.14
$ws.Range("N141").Value = -225290.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1039.2142
$ws.Range("I107").Value = 549.93335
$ws.Range("J107").Value = 1603.7693
$ws.Range("K107").Value = 1649.80005
$ws.Range("L107").Value = 4811.3079
$ws.Range("M107").Value = 270.1999499999999
$ws.Range("N107").Value = -8651.3079
$ws.Range("H140").Value = 18519512
$ws.Range("I140").Value = 26316528
$ws.Range("J140").Value = 1593.75
$ws.Range("K140").Value = 78949584
$ws.Range("L140").Value = 4781.25
$ws.Range("M140").Value = -78944404
$ws.Range("N140").Value = -15141.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1630.3478
$ws.Range("I97").Value = 1577.7778
$ws.Range("K97").Value = 1577.7778
$ws.Range("M97").Value = -1081.7778
$ws.Range("H113").Value = 4537.6
$ws.Range("I113").Value = 4359.5
$ws.Range("J113").Value = 5250
$ws.Range("K113").Value = 4359.5
$ws.Range("L113").Value = 5250
$ws.Range("M113").Value = -2189.5
$ws.Range("N113").Value = -9590

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 8000
$ws.Range("J4").Value = 8000
$ws.Range("L4").Value = 8000
$ws.Range("N4").Value = -8226
$ws.Range("H28").Value = 8000
$ws.Range("J28").Value = 8000
$ws.Range("L28").Value = 8000
$ws.Range("N28").Value = -8464
$ws.Range("H37").Value = 8000
$ws.Range("J37").Value = 8000
$ws.Range("L37").Value = 8000
$ws.Range("N37").Value = -8214
$ws.Range("H40").Value = 25006728
$ws.Range("I40").Value = 31256880
$ws.Range("K40").Value = 31256880
$ws.Range("M40").Value = -31256744
$ws.Range("H61").Value = 5692.737
$ws.Range("I61").Value = 5891.9414
$ws.Range("K61").Value = 5891.9414
$ws.Range("M61").Value = -5689.9414
$ws.Range("H69").Value = 40666.332
$ws.Range("J69").Value = 40666.332
$ws.Range("L69").Value = 40666.332
$ws.Range("N69").Value = -42288.332
$ws.Range("H72").Value = 40666.332
$ws.Range("J72").Value = 40666.332
$ws.Range("L72").Value = 121998.996
$ws.Range("N72").Value = -130110.996
$ws.Range("H113").Value = 5692.737
$ws.Range("I113").Value = 5891.9414
$ws.Range("K113").Value = 5891.9414
$ws.Range("M113").Value = -3721.9414
$ws.Range("H122").Value = 50006150
$ws.Range("I122").Value = 71434800
$ws.Range("J122").Value = 5962.6665
$ws.Range("K122").Value = 214304400
$ws.Range("L122").Value = 17887.9995
$ws.Range("M122").Value = -214301950
$ws.Range("N122").Value = -22787.9995
$ws.Range("H132").Value = 6267.347
$ws.Range("J132").Value = 7237.3794
$ws.Range("L132").Value = 21712.1382
$ws.Range("N132").Value = -26772.1382

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12051.833
$ws.Range("J81").Value = 15163.725
$ws.Range("L81").Value = 30327.45
$ws.Range("N81").Value = -32449.45
$ws.Range("H84").Value = 12051.833
$ws.Range("J84").Value = 15163.725
$ws.Range("L84").Value = 151637.25
$ws.Range("N84").Value = -162245.25
$ws.Range("H105").Value = 65000
$ws.Range("J105").Value = 65000
$ws.Range("L105").Value = 65000
$ws.Range("N105").Value = -71988
$ws.Range("H113").Value = 1619.2307
$ws.Range("J113").Value = 1174
$ws.Range("L113").Value = 3522
$ws.Range("N113").Value = -7862
$ws.Range("H122").Value = 3271.68
$ws.Range("I122").Value = 3471.375
$ws.Range("K122").Value = 10414.125
$ws.Range("M122").Value = -7964.125
$ws.Range("H132").Value = 2557.8223
$ws.Range("I132").Value = 1972.4324
$ws.Range("J132").Value = 5265.25
$ws.Range("K132").Value = 5917.2972
$ws.Range("L132").Value = 15795.75
$ws.Range("M132").Value = -3387.2972
$ws.Range("N132").Value = -20855.75
